$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.979.95"
$ws.Range("E2").Value = "  +4.56%  "

$ws.Range("D3").Value = "3.530.06"
$ws.Range("E3").Value = "  +4.92%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "188.72"
$ws.Range("E5").Value = "  +8.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "563.02"
$ws.Range("E6").Value = "  +7.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("E7").Value = "  +5.74%  "

$ws.Range("D8").Value = "3.520.64"
$ws.Range("E8").Value = "  +4.68%  "

$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.634"
$ws.Range("E10").Value = "  +4.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.156"
$ws.Range("E11").Value = "  +15.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.79"
$ws.Range("E12").Value = "  +2.62%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000275"
$ws.Range("E13").Value = "  +7.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.36"
$ws.Range("E14").Value = "  +3.03%  "

$ws.Range("D15").Value = "4.095.27"
$ws.Range("E15").Value = "  +5.00%  "

$ws.Range("D16").Value = "3.535.50"
$ws.Range("E16").Value = "  +5.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.63"
$ws.Range("E17").Value = "  +6.01%  "

$ws.Range("E18").Value = "  +3.10%  "

$ws.Range("D19").Value = "67.016.69"
$ws.Range("E19").Value = "  +4.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.09"
$ws.Range("E20").Value = "  +7.18%  "

$ws.Range("E21").Value = "  +3.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "421.38"
$ws.Range("E22").Value = "  +12.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.10"
$ws.Range("E23").Value = "  +10.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.78"
$ws.Range("E24").Value = "  +5.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.17"
$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.05"
$ws.Range("E26").Value = "  -4.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.90"
$ws.Range("E27").Value = "  +7.44%  "

$ws.Range("E28").Value = "  +8.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.11"
$ws.Range("E29").Value = "  -0.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.11"
$ws.Range("E30").Value = "  +10.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.40"
$ws.Range("E31").Value = "  +5.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "637.03"
$ws.Range("E32").Value = "  +0.67%  "

$ws.Range("E33").Value = "  +3.29%  "

$ws.Range("E34").Value = "  +4.63%  "

$ws.Range("E35").Value = "  +5.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "60.32"
$ws.Range("E36").Value = "  +3.88%  "

$ws.Range("D37").Value = "0.0₃0827"
$ws.Range("E37").Value = "  +11.21%  "

$ws.Range("E38").Value = "  +18.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.37"
$ws.Range("E39").Value = "  +5.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.387"
$ws.Range("E41").Value = "  +1.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.38"
$ws.Range("E42").Value = "  +13.58%  "

$ws.Range("D43").Value = "3.129.14"
$ws.Range("E43").Value = "  +4.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("E45").Value = "  -0.94%  "

$ws.Range("E46").Value = "  +9.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.38"
$ws.Range("E47").Value = "  +11.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0419"
$ws.Range("E48").Value = "  +5.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.75"
$ws.Range("E49").Value = "  +2.38%  "

$ws.Range("E50").Value = "  +5.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "140.24"
$ws.Range("E51").Value = "  +2.24%  "
